$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1683.2
$ws.Range("J40").Value = 1966.3334
$ws.Range("L40").Value = 1966.3334
$ws.Range("N40").Value = -2316.3334
$ws.Range("H113").Value = 3476.25
$ws.Range("J113").Value = 4450
$ws.Range("L113").Value = 4450
$ws.Range("N113").Value = -10958
$ws.Range("H138").Value = 5369.3477
$ws.Range("J138").Value = 5911.7646
$ws.Range("L138").Value = 17735.2938
$ws.Range("N138").Value = -28015.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1888.9783
$ws.Range("I32").Value = 1886.2
$ws.Range("J32").Value = 2014
$ws.Range("K32").Value = 1886.2
$ws.Range("L32").Value = 2014
$ws.Range("M32").Value = -1599.2
$ws.Range("N32").Value = -2588
$ws.Range("H45").Value = 4654.625
$ws.Range("I45").Value = 4654.625
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4654.625
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4277.625
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 3664.3333
$ws.Range("J46").Value = 2998
$ws.Range("L46").Value = 2998
$ws.Range("N46").Value = -3636
$ws.Range("H110").Value = 2802.1765
$ws.Range("I110").Value = 2922.5334
$ws.Range("K110").Value = 2922.5334
$ws.Range("M110").Value = -877.5333999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 975.3333
$ws.Range("I64").Value = 983
$ws.Range("J64").Value = 960
$ws.Range("K64").Value = 983
$ws.Range("L64").Value = 960
$ws.Range("M64").Value = -758
$ws.Range("N64").Value = -1410
$ws.Range("H67").Value = 975.3333
$ws.Range("I67").Value = 983
$ws.Range("J67").Value = 960
$ws.Range("K67").Value = 983
$ws.Range("L67").Value = 960
$ws.Range("M67").Value = -203
$ws.Range("N67").Value = -2520
$ws.Range("H99").Value = 1922.5
$ws.Range("I99").Value = 1930
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 1930
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -432
$ws.Range("N99").Value = -4896
$ws.Range("H107").Value = 4943
$ws.Range("I107").Value = 4941.6
$ws.Range("K107").Value = 4941.6
$ws.Range("M107").Value = -3021.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3027.2856
$ws.Range("I31").Value = 3599.2222
$ws.Range("K31").Value = 3599.2222
$ws.Range("M31").Value = -3304.2222
$ws.Range("H34").Value = 3027.2856
$ws.Range("I34").Value = 3599.2222
$ws.Range("K34").Value = 3599.2222
$ws.Range("M34").Value = -3397.2222
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 200390
$ws.Range("J100").Value = 200390
$ws.Range("L100").Value = 200390
$ws.Range("N100").Value = -202554
$ws.Range("H107").Value = 1669.3
$ws.Range("I107").Value = 1699.4286
$ws.Range("K107").Value = 1699.4286
$ws.Range("M107").Value = 220.5714
$ws.Range("H122").Value = 5145.25
$ws.Range("I122").Value = 3527.3333
$ws.Range("K122").Value = 10581.9999
$ws.Range("M122").Value = -8131.999899999999
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3504.25
$ws.Range("I132").Value = 3742.9
$ws.Range("K132").Value = 11228.7
$ws.Range("M132").Value = -8698.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1241.8
$ws.Range("J5").Value = 1324.2222
$ws.Range("L5").Value = 3972.6666
$ws.Range("N5").Value = -4196.6666
$ws.Range("H80").Value = 989.5
$ws.Range("I80").Value = 989.5
$ws.Range("K80").Value = 2968.5
$ws.Range("M80").Value = -2032.5
$ws.Range("H83").Value = 989.5
$ws.Range("I83").Value = 989.5
$ws.Range("K83").Value = 8905.5
$ws.Range("M83").Value = -4225.5
$ws.Range("H122").Value = 330.85
$ws.Range("J122").Value = 485.45456
$ws.Range("L122").Value = 4369.09104
$ws.Range("N122").Value = -9269.091039999999
$ws.Range("H127").Value = 1989.5
$ws.Range("J127").Value = 1989.5
$ws.Range("L127").Value = 5968.5
$ws.Range("N127").Value = -15888.5
$ws.Range("H135").Value = 1241.8
$ws.Range("J135").Value = 1324.2222
$ws.Range("L135").Value = 11917.9998
$ws.Range("N135").Value = -16987.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H102").Value = 3957
$ws.Range("I102").Value = 3821.25
$ws.Range("K102").Value = 3821.25
$ws.Range("M102").Value = -2199.25
$ws.Range("H107").Value = 1703.52
$ws.Range("I107").Value = 1171.7142
$ws.Range("K107").Value = 1171.7142
$ws.Range("M107").Value = 748.2858000000001
$ws.Range("H122").Value = 3688.7222
$ws.Range("I122").Value = 3950.1
$ws.Range("J122").Value = 3362
$ws.Range("K122").Value = 11850.3
$ws.Range("L122").Value = 10086
$ws.Range("M122").Value = -9400.299999999999
$ws.Range("N122").Value = -14986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8317.5
$ws.Range("J7").Value = 8827.143
$ws.Range("L7").Value = 8827.143
$ws.Range("N7").Value = -9051.143
$ws.Range("H40").Value = 5798.5
$ws.Range("I40").Value = 5698.25
$ws.Range("K40").Value = 5698.25
$ws.Range("M40").Value = -5562.25
$ws.Range("H46").Value = 127875
$ws.Range("I46").Value = 251250
$ws.Range("K46").Value = 251250
$ws.Range("M46").Value = -251062
$ws.Range("H55").Value = 228
$ws.Range("I55").Value = 149.5
$ws.Range("J55").Value = 280.33334
$ws.Range("K55").Value = 149.5
$ws.Range("L55").Value = 280.33334
$ws.Range("M55").Value = 23.5
$ws.Range("N55").Value = -626.33334
$ws.Range("H68").Value = 2939.9333
$ws.Range("I68").Value = 2971.4285
$ws.Range("J68").Value = 2912.375
$ws.Range("K68").Value = 2971.4285
$ws.Range("L68").Value = 2912.375
$ws.Range("M68").Value = -2222.4285
$ws.Range("N68").Value = -4410.375
$ws.Range("H71").Value = 2939.9333
$ws.Range("I71").Value = 2971.4285
$ws.Range("J71").Value = 2912.375
$ws.Range("K71").Value = 14857.1425
$ws.Range("L71").Value = 14561.875
$ws.Range("M71").Value = -11113.1425
$ws.Range("N71").Value = -22049.875
$ws.Range("H126").Value = 8317.5
$ws.Range("J126").Value = 8827.143
$ws.Range("L126").Value = 26481.429
$ws.Range("N126").Value = -31421.429
$ws.Range("H132").Value = 4002.348
$ws.Range("I132").Value = 4347.5
$ws.Range("K132").Value = 13042.5
$ws.Range("M132").Value = -10512.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2355.7144
$ws.Range("I122").Value = 2355.7144
$ws.Range("K122").Value = 7067.1432
$ws.Range("M122").Value = -4617.1432
$ws.Range("H126").Value = 1794.1875
$ws.Range("I126").Value = 1840.4667
$ws.Range("K126").Value = 5521.4001
$ws.Range("M126").Value = -3051.4001
$ws.Range("H132").Value = 5660
$ws.Range("I132").Value = 6097.0435
$ws.Range("J132").Value = 3147
$ws.Range("K132").Value = 18291.1305
$ws.Range("L132").Value = 9441
$ws.Range("M132").Value = -15761.1305
$ws.Range("N132").Value = -14501
$ws.Range("H136").Value = 6013.2856
$ws.Range("I136").Value = 6914.5
$ws.Range("J136").Value = 3129.4
$ws.Range("K136").Value = 20743.5
$ws.Range("L136").Value = 9388.200000000001
$ws.Range("M136").Value = -18193.5
$ws.Range("N136").Value = -14488.2
